# Apply the changes described by the commit:
#  1. Rename the "USB-UART" component to "Arduino-UART" (Sheet1!A5).
#  2. Move the active selection on Sheet1 from D17 to C5.
#  3. Turn off "print cell comments" (cellComments atEnd -> none) on every
#     sheet's page setup.

$wb = $excel.ActiveWorkbook

# --- 1. Update the component name ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5").Value = "Arduino-UART"

# --- 2. Fix up the selection on Sheet1 -------------------------------------------
$ws1.Range("C5").Select()

# --- 3. Stop printing comments on every sheet ------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ps = $ws.PageSetup
    # Re-assert the page setup values that are normally preserved so only
    # PrintComments actually changes in the rendered page setup.
    $ps.PaperSize = 9
    $ps.Zoom = 100
    $ps.FitToPagesWide = 1
    $ps.FitToPagesTall = 1
    $ps.Orientation = 1
    $ps.PrintComments = 0
}
